$wb = $excel.ActiveWorkbook

$sheetData = @{
    "ALC" = @{
        "H19" = 8403697
        "J19" = 465.5
        "L19" = 465.5
        "N19" = -815.5
        "H98" = 4641.3335
        "I98" = 4966.273
        "J98" = 3747.75
        "K98" = 4966.273
        "L98" = 3747.75
        "M98" = -3468.273
        "N98" = -6743.75
        "H116" = 3586.861
        "I116" = 3290.04
        "J116" = 4261.4546
        "K116" = 3290.04
        "L116" = 4261.4546
        "M116" = 151.96
        "N116" = -11145.4546
        "H122" = 4641.3335
        "I122" = 4966.273
        "J122" = 3747.75
        "K122" = 14898.819
        "L122" = 11243.25
        "M122" = -12448.819
        "N122" = -16143.25
        "H132" = 4446857.5
        "I132" = 5130350
        "K132" = 15391050
        "M132" = -15388520
        "H137" = 2225979.8
        "I137" = 2781308.2
        "K137" = 8343924.600000001
        "M137" = -8341374.600000001
        "H141" = 381354.06
        "I141" = 1788.9231
        "J141" = 564107.7
        "K141" = 5366.7693
        "L141" = 1692323.1
        "M141" = -186.7692999999999
        "N141" = -1702683.1
    }
    "ARM" = @{
        "H2" = 15627058
        "I2" = 19232072
        "J2" = 5333.3335
        "K2" = 19232072
        "L2" = 5333.3335
        "M2" = -19231959
        "N2" = -5559.3335
        "H32" = 6738.036
        "I32" = 5828.275
        "J32" = 30998.334
        "K32" = 5828.275
        "L32" = 30998.334
        "M32" = -5541.275
        "N32" = -31572.334
        "H34" = 40009.332
        "J34" = 40009.332
        "L34" = 40009.332
        "N34" = -40551.332
        "H45" = 1451.1774
        "I45" = 1118.7046
        "J45" = 2263.889
        "K45" = 1118.7046
        "L45" = 2263.889
        "M45" = -741.7046
        "N45" = -3017.889
        "H63" = 2476.4707
        "I63" = 2115.3845
        "J63" = 3650
        "K63" = 2115.3845
        "L63" = 3650
        "M63" = -1429.3845
        "N63" = -5022
        "H66" = 2476.4707
        "I66" = 2115.3845
        "J66" = 3650
        "K66" = 10576.9225
        "L66" = 18250
        "M66" = -7144.922500000001
        "N66" = -25114
        "H94" = 30030
        "J94" = 30030
        "L94" = 30030
        "N94" = -31832
        "H110" = 1248.7333
        "I110" = 514.5417
        "K110" = 514.5417
        "M110" = 1530.4583
        "H116" = 15627058
        "I116" = 19232072
        "J116" = 5333.3335
        "K116" = 19232072
        "L116" = 5333.3335
        "M116" = -19229778
        "N116" = -9921.333500000001
        "H132" = 2482.3333
        "I132" = 1885.3704
        "K132" = 5656.1112
        "M132" = -3126.1112
    }
    "BSM" = @{
        "H3" = 15627058
        "I3" = 19232072
        "J3" = 5333.3335
        "K3" = 19232072
        "L3" = 5333.3335
        "M3" = -19231958
        "N3" = -5561.3335
        "H105" = 1485.3784
        "I105" = 1477.95
        "J105" = 1494.1177
        "K105" = 1477.95
        "L105" = 1494.1177
        "M105" = 269.05
        "N105" = -4988.1177
        "H134" = 2507.08
        "I134" = 1880.1364
        "K134" = 5640.4092
        "M134" = -3105.4092
    }
    "CRP" = @{
        "H16" = 1312.3334
        "I16" = 834.1875
        "J16" = 2268.625
        "K16" = 834.1875
        "L16" = 2268.625
        "M16" = -547.1875
        "N16" = -2842.625
        "H22" = 1254.5454
        "I22" = 366.66666
        "K22" = 366.66666
        "M22" = -16.66665999999998
        "H31" = 1697160.2
        "I31" = 2224266.8
        "K31" = 2224266.8
        "M31" = -2223971.8
        "H34" = 1697160.2
        "I34" = 2224266.8
        "K34" = 2224266.8
        "M34" = -2224064.8
        "H58" = 14289313
        "I58" = 2618
        "J58" = 41672144
        "K58" = 2618
        "L58" = 41672144
        "M58" = -2415
        "N58" = -41672550
        "H99" = 3386
        "I99" = 1129.4286
        "J99" = 5642.5713
        "K99" = 1129.4286
        "L99" = 5642.5713
        "M99" = 368.5714
        "N99" = -8638.5713
        "H113" = 1312.3334
        "I113" = 834.1875
        "J113" = 2268.625
        "K113" = 834.1875
        "L113" = 2268.625
        "M113" = 1335.8125
        "N113" = -6608.625
        "H126" = 3386
        "I126" = 1129.4286
        "J126" = 5642.5713
        "K126" = 3388.2858
        "L126" = 16927.7139
        "M126" = -918.2857999999997
        "N126" = -21867.7139
        "H134" = 2554.8
        "I134" = 1435.4286
        "K134" = 4306.2858
        "M134" = -1771.2858
        "H136" = 14289313
        "I136" = 2618
        "J136" = 41672144
        "K136" = 7854
        "L136" = 125016432
        "M136" = -5304
        "N136" = -125021532
    }
    "CUL" = @{
        "H69" = 124052.664
        "I69" = 737
        "J69" = 159285.72
        "K69" = 2211
        "L69" = 477857.16
        "M69" = -1400
        "N69" = -479479.16
        "H72" = 124052.664
        "I72" = 737
        "J72" = 159285.72
        "K72" = 6633
        "L72" = 1433571.48
        "M72" = -2577
        "N72" = -1441683.48
    }
    "GSM" = @{
        "H80" = 4198.091
        "I80" = 3198.5
        "J80" = 4769.2856
        "K80" = 3198.5
        "L80" = 4769.2856
        "M80" = -2200.5
        "N80" = -6765.2856
        "H83" = 4198.091
        "I83" = 3198.5
        "J83" = 4769.2856
        "K83" = 15992.5
        "L83" = 23846.428
        "M83" = -11000.5
        "N83" = -33830.428
        "H113" = 3101.2222
        "I113" = 2670.3333
        "J113" = 3316.6667
        "K113" = 2670.3333
        "L113" = 3316.6667
        "M113" = -500.3332999999998
        "N113" = -7656.6667
        "H126" = 3433.238
        "I126" = 1296
        "J126" = 4501.857
        "K126" = 3888
        "L126" = 13505.571
        "M126" = -1418
        "N126" = -18445.571
    }
    "LTW" = @{
        "H46" = 1288.2693
        "I46" = 971.95123
        "J46" = 2467.2727
        "K46" = 971.95123
        "L46" = 2467.2727
        "M46" = -783.95123
        "N46" = -2843.2727
        "H132" = 5111.641
        "I132" = 1686.16
        "J132" = 11228.571
        "K132" = 5058.48
        "L132" = 33685.713
        "M132" = -2528.48
        "N132" = -38745.713
        "H136" = 2002039
        "I136" = 2565140.2
        "J136" = 5589.5454
        "K136" = 7695420.600000001
        "L136" = 16768.6362
        "M136" = -7692870.600000001
        "N136" = -21868.6362
    }
    "WVR" = @{
        "H107" = 1103.125
        "I107" = 450
        "J107" = 3933.3333
        "K107" = 1350
        "L107" = 11799.9999
        "M107" = 570
        "N107" = -15639.9999
        "H122" = 296125.28
        "I122" = 527887.9
        "J122" = 2559.3333
        "K122" = 1583663.7
        "L122" = 7677.999899999999
        "M122" = -1581213.7
        "N122" = -12577.9999
        "H132" = 120440.25
        "I132" = 143614.78
        "J132" = 30317.055
        "K132" = 430844.34
        "L132" = 90951.16500000001
        "M132" = -428314.34
        "N132" = -96011.16500000001
        "H136" = 1065.0889
        "I136" = 542.4722
        "J136" = 3155.5557
        "K136" = 1627.4166
        "L136" = 9466.667099999999
        "M136" = 922.5834
        "N136" = -14566.6671
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetData[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
